$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CompStat_1")

# --- Header: volume number 37 -> 38, and week dates 9/9-9/15 -> 9/16-9/22 ---
$ws.Range("A8").Value = "Volume 31   Number  38"
$ws.Range("C9").Value = "Report Covering the Week  9/16/2024  Through  9/22/2024"

# --- Column E got a touch wider to fit the new "***.*" / percentage text ---
$ws.Columns.Item(5).ColumnWidth = 7.433768

# --- Weekly Crime Complaints table (rows 14-31, 33): refreshed counts + recomputed % changes ---
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = "0"
$ws.Range("E14").Value = "***.*"
$ws.Range("F14").Value = 4
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = 300
$ws.Range("I14").Value = 13
$ws.Range("J14").Value = 16
$ws.Range("K14").Value = -18.75
$ws.Range("L14").Value = -35
$ws.Range("M14").Value = -35
$ws.Range("N14").Value = -86.734693877551
$ws.Range("C15").Value = 3
$ws.Range("D15").Value = 8
$ws.Range("E15").Value = -62.5
$ws.Range("F15").Value = 18
$ws.Range("G15").Value = 19
$ws.Range("H15").Value = -5.263157894736
$ws.Range("I15").Value = 160
$ws.Range("J15").Value = 146
$ws.Range("K15").Value = 9.589041095890
$ws.Range("L15").Value = 7.382550335570
$ws.Range("M15").Value = 70.212765957446
$ws.Range("N15").Value = 8.843537414965
$ws.Range("C16").Value = 53
$ws.Range("D16").Value = 53
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 167
$ws.Range("G16").Value = 197
$ws.Range("H16").Value = -15.228426395939
$ws.Range("I16").Value = 1608
$ws.Range("J16").Value = 1473
$ws.Range("K16").Value = 9.164969450101
$ws.Range("L16").Value = 22.935779816513
$ws.Range("M16").Value = 18.322295805739
$ws.Range("N16").Value = -74.780426599749
$ws.Range("C17").Value = 66
$ws.Range("D17").Value = 58
$ws.Range("E17").Value = 13.793103448275
$ws.Range("F17").Value = 235
$ws.Range("G17").Value = 243
$ws.Range("H17").Value = -3.292181069958
$ws.Range("I17").Value = 2444
$ws.Range("J17").Value = 2114
$ws.Range("K17").Value = 15.610217596972
$ws.Range("L17").Value = 33.917808219178
$ws.Range("M17").Value = 115.520282186949
$ws.Range("N17").Value = 13.200555812876
$ws.Range("C18").Value = 37
$ws.Range("D18").Value = 32
$ws.Range("E18").Value = 15.625
$ws.Range("F18").Value = 142
$ws.Range("G18").Value = 136
$ws.Range("H18").Value = 4.411764705882
$ws.Range("I18").Value = 1419
$ws.Range("J18").Value = 1424
$ws.Range("K18").Value = -0.351123595505
$ws.Range("L18").Value = 1.284796573875
$ws.Range("M18").Value = -25.628930817610
$ws.Range("N18").Value = -86.825735771980
$ws.Range("C19").Value = 123
$ws.Range("D19").Value = 127
$ws.Range("E19").Value = -3.149606299212
$ws.Range("F19").Value = 484
$ws.Range("G19").Value = 506
$ws.Range("H19").Value = -4.347826086956
$ws.Range("I19").Value = 4849
$ws.Range("J19").Value = 4967
$ws.Range("K19").Value = -2.375679484598
$ws.Range("L19").Value = -4.264560710760
$ws.Range("M19").Value = 68.426536992011
$ws.Range("N19").Value = -19.264069264069
$ws.Range("C20").Value = 58
$ws.Range("D20").Value = 88
$ws.Range("E20").Value = -34.090909090909
$ws.Range("F20").Value = 235
$ws.Range("G20").Value = 276
$ws.Range("H20").Value = -14.855072463768
$ws.Range("I20").Value = 1935
$ws.Range("J20").Value = 1910
$ws.Range("K20").Value = 1.308900523560
$ws.Range("L20").Value = 44.295302013422
$ws.Range("M20").Value = 48.389570552147
$ws.Range("N20").Value = -88.654353562005
$ws.Range("C21").Value = 341
$ws.Range("D21").Value = 366
$ws.Range("E21").Value = -6.830601092896
$ws.Range("F21").Value = 1285
$ws.Range("G21").Value = 1378
$ws.Range("H21").Value = -6.748911465892
$ws.Range("I21").Value = 12428
$ws.Range("J21").Value = 12050
$ws.Range("K21").Value = 3.136929460580
$ws.Range("L21").Value = 11.873255918624
$ws.Range("M21").Value = 42.883421476201
$ws.Range("N21").Value = -70.834506711724
$ws.Range("C22").Value = 7
$ws.Range("D22").Value = 10
$ws.Range("E22").Value = -30
$ws.Range("F22").Value = 21
$ws.Range("G22").Value = 29
$ws.Range("H22").Value = -27.586206896551
$ws.Range("I22").Value = 227
$ws.Range("J22").Value = 266
$ws.Range("K22").Value = -14.661654135338
$ws.Range("L22").Value = 7.075471698113
$ws.Range("M22").Value = 60.992907801418
$ws.Range("N22").Value = "***.*"
$ws.Range("C23").Value = 6
$ws.Range("D23").Value = 4
$ws.Range("E23").Value = 50
$ws.Range("F23").Value = 19
$ws.Range("G23").Value = 23
$ws.Range("H23").Value = -17.391304347826
$ws.Range("I23").Value = 188
$ws.Range("J23").Value = 202
$ws.Range("K23").Value = -6.930693069306
$ws.Range("L23").Value = 6.214689265536
$ws.Range("M23").Value = 51.612903225806
$ws.Range("N23").Value = "***.*"
$ws.Range("C24").Value = 262
$ws.Range("D24").Value = 318
$ws.Range("E24").Value = -17.610062893081
$ws.Range("F24").Value = 1165
$ws.Range("G24").Value = 1254
$ws.Range("H24").Value = -7.097288676236
$ws.Range("I24").Value = 11719
$ws.Range("J24").Value = 11300
$ws.Range("K24").Value = 3.707964601769
$ws.Range("L24").Value = 5.805344889851
$ws.Range("M24").Value = 71.933685446009
$ws.Range("N24").Value = "***.*"
$ws.Range("C25").Value = 160
$ws.Range("D25").Value = 187
$ws.Range("E25").Value = -14.438502673796
$ws.Range("F25").Value = 722
$ws.Range("G25").Value = 695
$ws.Range("H25").Value = 3.884892086330
$ws.Range("I25").Value = 7239
$ws.Range("J25").Value = 6299
$ws.Range("K25").Value = 14.923003651373
$ws.Range("L25").Value = 27.245561610124
$ws.Range("M25").Value = "***.*"
$ws.Range("N25").Value = "***.*"
$ws.Range("C26").Value = 121
$ws.Range("D26").Value = 127
$ws.Range("E26").Value = -4.724409448818
$ws.Range("F26").Value = 487
$ws.Range("G26").Value = 508
$ws.Range("H26").Value = -4.133858267716
$ws.Range("I26").Value = 4638
$ws.Range("J26").Value = 3902
$ws.Range("K26").Value = 18.862121988723
$ws.Range("L26").Value = 29.408482142857
$ws.Range("M26").Value = 31.911262798634
$ws.Range("N26").Value = "***.*"
$ws.Range("C27").Value = 5
$ws.Range("D27").Value = 10
$ws.Range("E27").Value = -50
$ws.Range("F27").Value = 25
$ws.Range("G27").Value = 28
$ws.Range("H27").Value = -10.714285714285
$ws.Range("I27").Value = 235
$ws.Range("J27").Value = 229
$ws.Range("K27").Value = 2.620087336244
$ws.Range("L27").Value = 9.302325581395
$ws.Range("M27").Value = "***.*"
$ws.Range("N27").Value = "***.*"
$ws.Range("C28").Value = 21
$ws.Range("D28").Value = 6
$ws.Range("E28").Value = 250
$ws.Range("F28").Value = 61
$ws.Range("G28").Value = 57
$ws.Range("H28").Value = 7.017543859649
$ws.Range("I28").Value = 462
$ws.Range("J28").Value = 530
$ws.Range("K28").Value = -12.830188679245
$ws.Range("L28").Value = 1.094091903719
$ws.Range("M28").Value = "***.*"
$ws.Range("N28").Value = "***.*"
$ws.Range("C29").Value = 1
$ws.Range("D29").Value = "0"
$ws.Range("E29").Value = "***.*"
$ws.Range("F29").Value = 4
$ws.Range("G29").Value = 4
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 20
$ws.Range("J29").Value = 50
$ws.Range("K29").Value = -60
$ws.Range("L29").Value = -66.101694915254
$ws.Range("M29").Value = -44.444444444444
$ws.Range("N29").Value = -89.583333333333
$ws.Range("C30").Value = 1
$ws.Range("D30").Value = "0"
$ws.Range("E30").Value = "***.*"
$ws.Range("F30").Value = 4
$ws.Range("G30").Value = 4
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 17
$ws.Range("J30").Value = 45
$ws.Range("K30").Value = -62.222222222222
$ws.Range("L30").Value = -66
$ws.Range("M30").Value = -45.161290322580
$ws.Range("N30").Value = -90.229885057471
$ws.Range("C31").Value = "0"
$ws.Range("D31").Value = 3
$ws.Range("E31").Value = -100
$ws.Range("F31").Value = 1
$ws.Range("G31").Value = 7
$ws.Range("H31").Value = -85.714285714285
$ws.Range("I31").Value = 43
$ws.Range("J31").Value = 59
$ws.Range("K31").Value = -27.118644067796
$ws.Range("L31").Value = -14
$ws.Range("M31").Value = "***.*"
$ws.Range("N31").Value = "***.*"
$ws.Range("C33").Value = "0"
$ws.Range("D33").Value = "0"
$ws.Range("E33").Value = "***.*"
$ws.Range("F33").Value = 4
$ws.Range("G33").Value = 1
$ws.Range("H33").Value = 300
$ws.Range("I33").Value = 36
$ws.Range("J33").Value = 34
$ws.Range("K33").Value = 5.882352941176
$ws.Range("L33").Value = 38.461538461538
$ws.Range("M33").Value = "***.*"
$ws.Range("N33").Value = "***.*"
